$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1453.091
$ws.Range("J17").Value = 1503.4
$ws.Range("L17").Value = 4510.200000000001
$ws.Range("N17").Value = -4846.200000000001
$ws.Range("H40").Value = 4260.727
$ws.Range("I40").Value = 3294.8333
$ws.Range("K40").Value = 3294.8333
$ws.Range("M40").Value = -3119.8333
$ws.Range("H69").Value = 53338.332
$ws.Range("I69").Value = 20000
$ws.Range("J69").Value = 70007.5
$ws.Range("K69").Value = 60000
$ws.Range("L69").Value = 210022.5
$ws.Range("M69").Value = -59126
$ws.Range("N69").Value = -211770.5
$ws.Range("H72").Value = 53338.332
$ws.Range("I72").Value = 20000
$ws.Range("J72").Value = 70007.5
$ws.Range("K72").Value = 180000
$ws.Range("L72").Value = 630067.5
$ws.Range("M72").Value = -175632
$ws.Range("N72").Value = -638803.5
$ws.Range("H76").Value = 3934.889
$ws.Range("I76").Value = 3934.889
$ws.Range("K76").Value = 3934.889
$ws.Range("M76").Value = -3619.889
$ws.Range("H79").Value = 3934.889
$ws.Range("I79").Value = 3934.889
$ws.Range("K79").Value = 3934.889
$ws.Range("M79").Value = -2842.889
$ws.Range("H111").Value = 2132.125
$ws.Range("I111").Value = 1676.1666
$ws.Range("K111").Value = 5028.4998
$ws.Range("M111").Value = -1961.4998
$ws.Range("H129").Value = 2369.8333
$ws.Range("J129").Value = 3205.6667
$ws.Range("L129").Value = 9617.000100000001
$ws.Range("N129").Value = -19617.0001
$ws.Range("H132").Value = 2940.9429
$ws.Range("J132").Value = 1837.5
$ws.Range("L132").Value = 5512.5
$ws.Range("N132").Value = -10572.5
$ws.Range("H135").Value = 1253
$ws.Range("I135").Value = 1253
$ws.Range("K135").Value = 11277
$ws.Range("M135").Value = -8742
$ws.Range("H137").Value = 3272.4255
$ws.Range("I137").Value = 2146.5667
$ws.Range("K137").Value = 6439.7001
$ws.Range("M137").Value = -3889.7001
$ws.Range("H138").Value = 3443.432
$ws.Range("I138").Value = 1486.5
$ws.Range("J138").Value = 3639.125
$ws.Range("K138").Value = 4459.5
$ws.Range("L138").Value = 10917.375
$ws.Range("M138").Value = 680.5
$ws.Range("N138").Value = -21197.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10640889
$ws.Range("I32").Value = 12196752
$ws.Range("K32").Value = 12196752
$ws.Range("M32").Value = -12196465
$ws.Range("H41").Value = 5625
$ws.Range("I41").Value = 5625
$ws.Range("K41").Value = 5625
$ws.Range("M41").Value = -5211
$ws.Range("H61").Value = 17897588
$ws.Range("I61").Value = 20003522
$ws.Range("K61").Value = 20003522
$ws.Range("M61").Value = -20003310
$ws.Range("H75").Value = 36500
$ws.Range("J75").Value = 36500
$ws.Range("L75").Value = 36500
$ws.Range("N75").Value = -38248
$ws.Range("H78").Value = 36500
$ws.Range("J78").Value = 36500
$ws.Range("L78").Value = 109500
$ws.Range("N78").Value = -118236
$ws.Range("H136").Value = 17897588
$ws.Range("I136").Value = 20003522
$ws.Range("K136").Value = 60010566
$ws.Range("M136").Value = -60008016

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 511.81818
$ws.Range("I94").Value = 329.125
$ws.Range("K94").Value = 329.125
$ws.Range("M94").Value = 121.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 773.6429000000001
$ws.Range("I16").Value = 646.36365
$ws.Range("K16").Value = 646.36365
$ws.Range("M16").Value = -359.36365
$ws.Range("H113").Value = 773.6429000000001
$ws.Range("I113").Value = 646.36365
$ws.Range("K113").Value = 646.36365
$ws.Range("M113").Value = 1523.63635
$ws.Range("H122").Value = 1108.25
$ws.Range("I122").Value = 1108.25
$ws.Range("K122").Value = 3324.75
$ws.Range("M122").Value = -874.75
$ws.Range("H132").Value = 1918.0968
$ws.Range("I132").Value = 1659.3572
$ws.Range("K132").Value = 4978.071599999999
$ws.Range("M132").Value = -2448.071599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 133.1923
$ws.Range("I2").Value = 108.71429
$ws.Range("J2").Value = 142.21053
$ws.Range("K2").Value = 652.28574
$ws.Range("L2").Value = 853.26318
$ws.Range("M2").Value = -539.28574
$ws.Range("N2").Value = -1079.26318
$ws.Range("H4").Value = 13833574
$ws.Range("I4").Value = 20416862
$ws.Range("K4").Value = 61250586
$ws.Range("M4").Value = -61250474
$ws.Range("H23").Value = 361.86667
$ws.Range("I23").Value = 688.1429000000001
$ws.Range("J23").Value = 76.375
$ws.Range("K23").Value = 2064.4287
$ws.Range("L23").Value = 229.125
$ws.Range("M23").Value = -1829.4287
$ws.Range("N23").Value = -699.125
$ws.Range("H39").Value = 1800
$ws.Range("I39").Value = 1800
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 5400
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -5106
$ws.Range("H44").Value = 21760.75
$ws.Range("I44").Value = 24442.857
$ws.Range("K44").Value = 73328.571
$ws.Range("M44").Value = -72930.571
$ws.Range("H55").Value = 3868.3333
$ws.Range("I55").Value = 800
$ws.Range("K55").Value = 2400
$ws.Range("M55").Value = -2223
$ws.Range("H68").Value = 2049.6667
$ws.Range("I68").Value = 1532.8889
$ws.Range("K68").Value = 4598.6667
$ws.Range("M68").Value = -3787.6667
$ws.Range("H71").Value = 2049.6667
$ws.Range("I71").Value = 1532.8889
$ws.Range("K71").Value = 13796.0001
$ws.Range("M71").Value = -9740.000099999999
$ws.Range("H116").Value = 350
$ws.Range("I116").Value = 350
$ws.Range("K116").Value = 1050
$ws.Range("M116").Value = 2392
$ws.Range("H131").Value = 4082.225
$ws.Range("I131").Value = 6021.2144
$ws.Range("J131").Value = 3670.9243
$ws.Range("K131").Value = 18063.6432
$ws.Range("L131").Value = 11012.7729
$ws.Range("M131").Value = -13023.6432
$ws.Range("N131").Value = -21092.7729

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3134652.2
$ws.Range("J11").Value = 2425762.2
$ws.Range("L11").Value = 2425762.2
$ws.Range("N11").Value = -2426040.2
$ws.Range("H18").Value = 20000
$ws.Range("J18").Value = 20000
$ws.Range("L18").Value = 20000
$ws.Range("N18").Value = -20586

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 28334.277
$ws.Range("I20").Value = 19999.916
$ws.Range("J20").Value = 45003
$ws.Range("K20").Value = 19999.916
$ws.Range("L20").Value = 45003
$ws.Range("M20").Value = -19773.916
$ws.Range("N20").Value = -45455
$ws.Range("H40").Value = 2953.8386
$ws.Range("I40").Value = 1988.3
$ws.Range("J40").Value = 4709.364
$ws.Range("K40").Value = 1988.3
$ws.Range("L40").Value = 4709.364
$ws.Range("M40").Value = -1852.3
$ws.Range("N40").Value = -4981.364
$ws.Range("H46").Value = 4166.6333
$ws.Range("I46").Value = 1939.0555
$ws.Range("J46").Value = 7508
$ws.Range("K46").Value = 1939.0555
$ws.Range("L46").Value = 7508
$ws.Range("M46").Value = -1751.0555
$ws.Range("N46").Value = -7884
$ws.Range("H82").Value = 1221.3889
$ws.Range("I82").Value = 1155.25
$ws.Range("J82").Value = 1353.6666
$ws.Range("K82").Value = 1155.25
$ws.Range("L82").Value = 1353.6666
$ws.Range("M82").Value = -794.25
$ws.Range("N82").Value = -2075.6666
$ws.Range("H85").Value = 1221.3889
$ws.Range("I85").Value = 1155.25
$ws.Range("J85").Value = 1353.6666
$ws.Range("K85").Value = 1155.25
$ws.Range("L85").Value = 1353.6666
$ws.Range("M85").Value = 92.75
$ws.Range("N85").Value = -3849.6666
$ws.Range("H132").Value = 235792.89
$ws.Range("I132").Value = 2626.2942
$ws.Range("J132").Value = 1116644.5
$ws.Range("K132").Value = 7878.882599999999
$ws.Range("L132").Value = 3349933.5
$ws.Range("M132").Value = -5348.882599999999
$ws.Range("N132").Value = -3354993.5
$ws.Range("H140").Value = 135377
$ws.Range("J140").Value = 135377
$ws.Range("L140").Value = 135377
$ws.Range("N140").Value = -145737

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H100").Value = 487.25
$ws.Range("I100").Value = 481.66666
$ws.Range("K100").Value = 963.33332
$ws.Range("M100").Value = -422.33332
$ws.Range("H107").Value = 15625525
$ws.Range("I107").Value = 25000492
$ws.Range("K107").Value = 75001476
$ws.Range("M107").Value = -74999556
$ws.Range("H113").Value = 562.5217
$ws.Range("J113").Value = 499.4
$ws.Range("L113").Value = 1498.2
$ws.Range("N113").Value = -5838.2
$ws.Range("H122").Value = 1399.439
$ws.Range("I122").Value = 1443.8857
$ws.Range("J122").Value = 1140.1666
$ws.Range("K122").Value = 4331.6571
$ws.Range("L122").Value = 3420.4998
$ws.Range("M122").Value = -1881.6571
$ws.Range("N122").Value = -8320.4998
$ws.Range("H132").Value = 1635.2858
$ws.Range("I132").Value = 1485.898
$ws.Range("J132").Value = 2681
$ws.Range("K132").Value = 4457.694
$ws.Range("L132").Value = 8043
$ws.Range("M132").Value = -1927.694
$ws.Range("N132").Value = -13103
